$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the target (budget) value
$ws.Range("H2").Value = 12000

# Update quantities in column G, which drive the formulas in column H
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 0

# Update the selected cell to match the diff
$ws.Range("G10").Select()
